# Generate Report for Handoff
# Re-applies the localization-status refresh: the "98a1ad4e-...md" file moves
# from "Handed back: in sync with en-US" back to "Ready for handoff" (new
# handoff timestamps), while "ffff0fa14419-...md" and "ffffff22ef7109-...md"
# keep "Handed back" status. Rows on every sheet are reordered so that
# 98a1ad4e now sits just above the ".localization-config" row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "ffff0fa14419-b49a-436c-bb8c-abf194bdc1e9.md"
$ws1.Range("B2").Value = "Handed back: in sync with en-US"
$ws1.Range("C2").Value = "Handed back: in sync with en-US"

$ws1.Range("A3").Value = "ffffff22ef7109-8d16-4492-af9a-6a08ffbe39eb.md"
$ws1.Range("B3").Value = "Handed back: in sync with en-US"
$ws1.Range("C3").Value = "Handed back: in sync with en-US"

$ws1.Range("A4").Value = "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.md"
$ws1.Range("B4").Value = "Ready for handoff"
$ws1.Range("C4").Value = "Ready for handoff"

$ws1.Range("A5").Value = ".localization-config"
$ws1.Range("B5").Value = "Not to be localized"
$ws1.Range("C5").Value = "Not to be localized"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7030d1bd6c5d3388755f1382eca3fbdc0438e01c/e2e/98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.md", "", "", "ffff0fa14419-b49a-436c-bb8c-abf194bdc1e9.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/7030d1bd6c5d3388755f1382eca3fbdc0438e01c/e2e/ffff0fa14419-b49a-436c-bb8c-abf194bdc1e9.md", "", "", "ffffff22ef7109-8d16-4492-af9a-6a08ffbe39eb.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/7030d1bd6c5d3388755f1382eca3fbdc0438e01c/e2e/ffffff22ef7109-8d16-4492-af9a-6a08ffbe39eb.md", "", "", "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7030d1bd6c5d3388755f1382eca3fbdc0438e01c/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "ffff0fa14419-b49a-436c-bb8c-abf194bdc1e9.md"
$ws2.Range("B2").Value = "Handed back: in sync with en-US"
$ws2.Range("C2").Value = "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-03-09 01:35:10"
$ws2.Range("E2").Value = "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.md"
$ws2.Range("F2").Value = "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.zh-cn.xlf"
$ws2.Range("G2").Value = "2016-03-09 01:36:22"
$ws2.Range("H2").Value = "Include"

$ws2.Range("A3").Value = "ffffff22ef7109-8d16-4492-af9a-6a08ffbe39eb.md"
$ws2.Range("B3").Value = "Handed back: in sync with en-US"
$ws2.Range("C3").Value = "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-03-09 01:35:10"
$ws2.Range("E3").Value = "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.md"
$ws2.Range("F3").Value = "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.zh-cn.xlf"
$ws2.Range("G3").Value = "2016-03-09 01:36:22"
$ws2.Range("H3").Value = "Include"

$ws2.Range("A4").Value = "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.md"
$ws2.Range("B4").Value = "Ready for handoff"
$ws2.Range("C4").Value = "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.28501e0d2cdabbe76c4070eb73eb9ef73e71bb4b.zh-cn.xlf"
$ws2.Range("D4").Value = "2016-03-09 01:40:31"
$ws2.Range("E4").Value = "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.md"
$ws2.Range("F4").Value = "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.28501e0d2cdabbe76c4070eb73eb9ef73e71bb4b.zh-cn.xlf"
$ws2.Range("G4").Value = "2016-03-09 01:39:08"
$ws2.Range("H4").Value = "Include"

$ws2.Range("A5").Value = ".localization-config"
$ws2.Range("B5").Value = "Not to be localized"
$ws2.Range("D5").Value = "0001-01-01 00:00:00"
$ws2.Range("G5").Value = "0001-01-01 00:00:00"
$ws2.Range("H5").Value = "Ignored"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7030d1bd6c5d3388755f1382eca3fbdc0438e01c/e2e/98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.md", "", "", "ffff0fa14419-b49a-436c-bb8c-abf194bdc1e9.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/faee87627a6f4a4fc3e368ad8bc634fc3ee8f1ab/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.28501e0d2cdabbe76c4070eb73eb9ef73e71bb4b.zh-cn.xlf", "", "", "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/a6771072a4b186fc1623624534a077a5c01df96d/e2e/98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.md", "", "", "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/0966c17051a61436664095b968e0c2204408f251/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.28501e0d2cdabbe76c4070eb73eb9ef73e71bb4b.zh-cn.xlf", "", "", "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/7030d1bd6c5d3388755f1382eca3fbdc0438e01c/e2e/ffff0fa14419-b49a-436c-bb8c-abf194bdc1e9.md", "", "", "ffffff22ef7109-8d16-4492-af9a-6a08ffbe39eb.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6613d1ac9a6a252ddcf64cd92770c9358894c03b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.zh-cn.xlf", "", "", "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/32644b487c547dec539f2b85275997a263b8e816/e2e/b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.md", "", "", "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7538f700d69db380b64977607f8b171a6b88264f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.zh-cn.xlf", "", "", "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/7030d1bd6c5d3388755f1382eca3fbdc0438e01c/e2e/ffffff22ef7109-8d16-4492-af9a-6a08ffbe39eb.md", "", "", "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6613d1ac9a6a252ddcf64cd92770c9358894c03b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.zh-cn.xlf", "", "", "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.28501e0d2cdabbe76c4070eb73eb9ef73e71bb4b.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/32644b487c547dec539f2b85275997a263b8e816/e2e/b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.md", "", "", "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7538f700d69db380b64977607f8b171a6b88264f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.zh-cn.xlf", "", "", "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.28501e0d2cdabbe76c4070eb73eb9ef73e71bb4b.zh-cn.xlf") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7030d1bd6c5d3388755f1382eca3fbdc0438e01c/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "ffff0fa14419-b49a-436c-bb8c-abf194bdc1e9.md"
$ws3.Range("B2").Value = "Handed back: in sync with en-US"
$ws3.Range("C2").Value = "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.de-de.xlf"
$ws3.Range("D2").Value = "2016-03-09 01:35:25"
$ws3.Range("E2").Value = "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.md"
$ws3.Range("F2").Value = "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.de-de.xlf"
$ws3.Range("G2").Value = "2016-03-09 01:36:53"
$ws3.Range("H2").Value = "Include"

$ws3.Range("A3").Value = "ffffff22ef7109-8d16-4492-af9a-6a08ffbe39eb.md"
$ws3.Range("B3").Value = "Handed back: in sync with en-US"
$ws3.Range("C3").Value = "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.de-de.xlf"
$ws3.Range("D3").Value = "2016-03-09 01:35:25"
$ws3.Range("E3").Value = "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.md"
$ws3.Range("F3").Value = "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.de-de.xlf"
$ws3.Range("G3").Value = "2016-03-09 01:36:53"
$ws3.Range("H3").Value = "Include"

$ws3.Range("A4").Value = "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.md"
$ws3.Range("B4").Value = "Ready for handoff"
$ws3.Range("C4").Value = "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.28501e0d2cdabbe76c4070eb73eb9ef73e71bb4b.de-de.xlf"
$ws3.Range("D4").Value = "2016-03-09 01:40:40"
$ws3.Range("E4").Value = "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.md"
$ws3.Range("F4").Value = "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.28501e0d2cdabbe76c4070eb73eb9ef73e71bb4b.de-de.xlf"
$ws3.Range("G4").Value = "2016-03-09 01:39:39"
$ws3.Range("H4").Value = "Include"

$ws3.Range("A5").Value = ".localization-config"
$ws3.Range("B5").Value = "Not to be localized"
$ws3.Range("D5").Value = "0001-01-01 00:00:00"
$ws3.Range("G5").Value = "0001-01-01 00:00:00"
$ws3.Range("H5").Value = "Ignored"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/7030d1bd6c5d3388755f1382eca3fbdc0438e01c/e2e/98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.md", "", "", "ffff0fa14419-b49a-436c-bb8c-abf194bdc1e9.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d987d00349b049a85359066b74d48a15f1b24548/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.28501e0d2cdabbe76c4070eb73eb9ef73e71bb4b.de-de.xlf", "", "", "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/d1a107d73fd10f0d7899b9364b21833790d07fdb/e2e/98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.md", "", "", "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/f302b569da683c810fc34a589fdc5406f4bebafe/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.28501e0d2cdabbe76c4070eb73eb9ef73e71bb4b.de-de.xlf", "", "", "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/7030d1bd6c5d3388755f1382eca3fbdc0438e01c/e2e/ffff0fa14419-b49a-436c-bb8c-abf194bdc1e9.md", "", "", "ffffff22ef7109-8d16-4492-af9a-6a08ffbe39eb.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f3eb10877f1b57fde5c98dc606fbd2b57c462a79/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.de-de.xlf", "", "", "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1389fc4e550d3844f3614d265b37ac3b894d0f50/e2e/b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.md", "", "", "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/63aaa824c08b72a598c8baa3547f20dc5aa77c58/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.de-de.xlf", "", "", "b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/7030d1bd6c5d3388755f1382eca3fbdc0438e01c/e2e/ffffff22ef7109-8d16-4492-af9a-6a08ffbe39eb.md", "", "", "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f3eb10877f1b57fde5c98dc606fbd2b57c462a79/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.de-de.xlf", "", "", "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.28501e0d2cdabbe76c4070eb73eb9ef73e71bb4b.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("E4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/1389fc4e550d3844f3614d265b37ac3b894d0f50/e2e/b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.md", "", "", "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("F4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/63aaa824c08b72a598c8baa3547f20dc5aa77c58/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b135e0ff-a1d1-42a4-ad07-0bb72b2b6910.ef3e64830d5546932623ab9473ece8ca9cd1230f.de-de.xlf", "", "", "98a1ad4e-1b0a-4c99-9ed4-1237322d6a1f.28501e0d2cdabbe76c4070eb73eb9ef73e71bb4b.de-de.xlf") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/7030d1bd6c5d3388755f1382eca3fbdc0438e01c/.localization-config", "", "", ".localization-config") | Out-Null

$wb.Save()
